$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 4).Value  = 0.382                 # D
    $ws.Cells.Item($r, 5).Value  = 0.457                 # E
    $ws.Cells.Item($r, 6).Value  = 0.217                 # F

    $ws.Cells.Item($r, 9).Value  = 0                     # I
    $ws.Cells.Item($r, 10).Value = 0                     # J
    $ws.Cells.Item($r, 11).Value = 19                    # K
    $ws.Cells.Item($r, 12).Value = 0.2592087312414734    # L
    $ws.Cells.Item($r, 13).Value = 12.4                  # M
    $ws.Cells.Item($r, 14).Value = 0.0625                # N
    $ws.Cells.Item($r, 15).Value = 0.6526315789473685    # O
    $ws.Cells.Item($r, 16).Value = 12.4                  # P
    $ws.Cells.Item($r, 17).Value = 0.0625                # Q
    $ws.Cells.Item($r, 18).Value = 0.6526315789473685    # R

    $ws.Cells.Item($r, 21).Value = 87                    # U
    $ws.Cells.Item($r, 22).Value = 0.438508064516129     # V
    $ws.Cells.Item($r, 23).Value = 0.236612702366127     # W
    $ws.Cells.Item($r, 24).Value = 0.05991588066491206   # X
    $ws.Cells.Item($r, 25).Value = 0.176696821701215     # Y
    $ws.Cells.Item($r, 26).Value = 7.881720430107529     # Z
    $ws.Cells.Item($r, 27).Value = 0                     # AA
    $ws.Cells.Item($r, 28).Value = 0.05907453693797746   # AB
    $ws.Cells.Item($r, 29).Value = -0.05907453693797746  # AC
    $ws.Cells.Item($r, 30).Value = 14.4                  # AD
    $ws.Cells.Item($r, 31).Value = 0                     # AE
    $ws.Cells.Item($r, 32).Value = 14.4                  # AF
    $ws.Cells.Item($r, 33).Value = -72.59999999999999    # AG
    $ws.Cells.Item($r, 34).Value = 0.06766917293233082   # AH
    $ws.Cells.Item($r, 35).Value = 0.1324747010119595    # AI
    $ws.Cells.Item($r, 36).Value = -0.5771065182829888   # AJ
    $ws.Cells.Item($r, 37).Value = -3.345622119815668    # AK

    # AN and AP are removed entirely (no longer present in the row)
    $ws.Cells.Item($r, 40).ClearContents()                # AN
    $ws.Cells.Item($r, 42).ClearContents()                # AP
}
